$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that contain the per-row "measurement" data that gets shifted down
# by one row (a new week's record is inserted at row 413, pushing older
# records down; the oldest record - formerly row 508 - becomes new row 509).
$shiftCols = @("D", "J", "K", "L", "M", "O", "P")

$firstRow = 413
$lastRow = 508
$newLastRow = 509

# 1) Read the current ("old") values for the columns that shift, for every
#    row from firstRow..lastRow, before we start overwriting anything.
$old = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $shiftCols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $old[$r] = $rowVals
}

# 2) Create the new last row (509) by copying the constant columns from the
#    previous last row (508) and the shifting columns from the old row 508.
$constCols = @("A", "B", "C", "E", "F", "G", "H", "I", "N", "Q", "R")
foreach ($c in $constCols) {
    $ws.Range("$c$newLastRow").Value = $ws.Range("$c$lastRow").Value()
}
foreach ($c in $shiftCols) {
    $ws.Range("$c$newLastRow").Value = $old[$lastRow][$c]
}
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastRow").NumberFormat

# 3) Shift rows 508 down to 414: row r gets the old values that used to be
#    in row (r-1).
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $prev = $r - 1
    foreach ($c in $shiftCols) {
        $ws.Range("$c$r").Value = $old[$prev][$c]
    }
}

# 4) Row 413 receives the brand new data point.
$ws.Range("D$firstRow").Value = 45204
$ws.Range("J$firstRow").Value = 500
$ws.Range("K$firstRow").Value = 1800
$ws.Range("L$firstRow").Value = 1800
$ws.Range("M$firstRow").Value = 1800
$ws.Range("O$firstRow").Value = "Región del Maule"
$ws.Range("P$firstRow").Value = 450
